$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Update text labels (shared-string edits from the diff)
$ws.Range("D2").Value = "Save Case File Rules"
$ws.Range("C12").Value = "RuleTable Save Case File Rules"

# Tighten row heights for rows 2 and 12 (ht 15 -> 13.3)
$ws.Rows.Item(2).RowHeight = 13.3
$ws.Rows.Item(12).RowHeight = 13.3

# Update the active selection / scroll position (topLeftCell A4->A1, selection C22->C13)
$ws.Range("C13").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
